# Daily attendance processing - 2025-12-08 07:37:56
# Updates "Recorded By" email orderings, session-count / coverage statistics,
# and flips the PATHOLOGY LAB/MUSEUM session-2 (08/12/2025) row from
# "Not Recorded" to "Recorded" with its attendee + attendance numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reordered "Recorded By" attendee lists -------------------------------
$ws.Range("G2").Value  = "gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G3").Value  = "hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G4").Value  = "hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G5").Value  = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G7").Value  = "menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G12").Value = "dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G27").Value = "nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
$ws.Range("G30").Value = "yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"

# --- ANATOMY class-statistics block (K2:L10) -------------------------------
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 2
$ws.Range("L9").Value = "79.3%"
$ws.Range("L10").Value = "27.1%"

# --- Matching per-group summary row for ANATOMY (row 15) -------------------
$ws.Range("O15").Value = 23
$ws.Range("P15").Value = 2
$ws.Range("R15").Value = "79.3%"
$ws.Range("S15").Value = "27.1%"

# --- PATHOLOGY LAB/MUSEUM session 2 (08/12/2025) newly recorded -----------
# Copy the "Recorded" (green) formatting from row 2 onto row 25 without
# touching row 25's own Year/Group/Subject/Session/Date/Time values.
$src = $ws.Range("A2:I2")
$dst = $ws.Range("A25:I25")
$src.Copy()
$dst.PasteSpecial(-4122)

$ws.Range("G25").Value = "menna-allah.gamil@med.asu.edu.eg"
$ws.Range("H25").Value = "60/251"
$ws.Range("I25").Value = "Recorded"
